$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: fill in previously-blank car record (Audi etron GT Quattro, row 21) ---
$ws.Range("C21").Value = 2021
$ws.Range("D21").Value = "4WD"
$ws.Range("E21").Value = 563
$ws.Range("F21").Value = "G"
$ws.Range("G21").Value = "légendaire"
$ws.Range("H21").Value = 0

# --- Row 33: fill in previously-blank car record (Chevrolet Camaro ZL1 1LE, row 33) ---
$ws.Range("D33").Value = "RWD"
$ws.Range("G33").Value = "épique"
$ws.Range("H33").Value = 0

# --- Row 87: fill in previously-blank car record (MercedesBenz 300 SL, row 87) ---
$ws.Range("D87").Value = "RWD"
$ws.Range("G87").Value = "légendaire"
$ws.Range("H87").Value = 0

# --- nb_vitesse (column I) values for every data row ---
$ws.Range("I3").Value = 4
$ws.Range("I9").Value = 5
$ws.Range("I10").Value = 7
$ws.Range("I11").Value = 6
$ws.Range("I12").Value = 8
$ws.Range("I13").Value = 8
$ws.Range("I14").Value = 6
$ws.Range("I15").Value = 7
$ws.Range("I16").Value = 7
$ws.Range("I17").Value = 7
$ws.Range("I18").Value = 7
$ws.Range("I19").Value = 7
$ws.Range("I20").Value = 7
$ws.Range("I21").Value = 2
$ws.Range("I22").Value = 6
$ws.Range("I23").Value = 8
$ws.Range("I24").Value = 6
$ws.Range("I25").Value = 8
$ws.Range("I26").Value = 7
$ws.Range("I27").Value = 7
$ws.Range("I28").Value = 7
$ws.Range("I29").Value = 7
$ws.Range("I30").Value = 7
$ws.Range("I31").Value = 6
$ws.Range("I32").Value = 6
$ws.Range("I33").Value = 6
$ws.Range("I34").Value = 3
$ws.Range("I35").Value = 7
$ws.Range("I36").Value = 7
$ws.Range("I37").Value = 4
$ws.Range("I38").Value = 6
$ws.Range("I39").Value = 6
$ws.Range("I40").Value = 5
$ws.Range("I41").Value = 5
$ws.Range("I42").Value = 5
$ws.Range("I43").Value = 7
$ws.Range("I44").Value = 7
$ws.Range("I45").Value = 6
$ws.Range("I46").Value = 5
$ws.Range("I47").Value = 6
$ws.Range("I48").Value = 7
$ws.Range("I49").Value = 7
$ws.Range("I50").Value = 10
$ws.Range("I51").Value = 4
$ws.Range("I52").Value = 6
$ws.Range("I53").Value = 6
$ws.Range("I54").Value = 6
$ws.Range("I55").Value = 6
$ws.Range("I56").Value = 7
$ws.Range("I57").Value = 7
$ws.Range("I58").Value = 7
$ws.Range("I59").Value = 7
$ws.Range("I60").Value = 5
$ws.Range("I61").Value = 4
$ws.Range("I62").Value = 4
$ws.Range("I63").Value = 8
$ws.Range("I64").Value = 7
$ws.Range("I65").Value = 1
$ws.Range("I68").Value = 7
$ws.Range("I69").Value = 7
$ws.Range("I70").Value = 5
$ws.Range("I71").Value = 7
$ws.Range("I72").Value = 7
$ws.Range("I73").Value = 5
$ws.Range("I74").Value = 8
$ws.Range("I75").Value = 8
$ws.Range("I76").Value = 6
$ws.Range("I77").Value = 1
$ws.Range("I78").Value = 6
$ws.Range("I79").Value = 6
$ws.Range("I80").Value = 8
$ws.Range("I81").Value = 7
$ws.Range("I82").Value = 7
$ws.Range("I83").Value = 7
$ws.Range("I84").Value = 8
$ws.Range("I85").Value = 7
$ws.Range("I86").Value = 7
$ws.Range("I87").Value = 4
$ws.Range("I88").Value = 7
$ws.Range("I89").Value = 5
$ws.Range("I90").Value = 7
$ws.Range("I91").Value = 6
$ws.Range("I92").Value = 6
$ws.Range("I93").Value = 5
$ws.Range("I94").Value = 6
$ws.Range("I95").Value = 7
$ws.Range("I96").Value = 2
$ws.Range("I97").Value = 8
$ws.Range("I98").Value = 4
$ws.Range("I99").Value = 6
$ws.Range("I100").Value = 4
$ws.Range("I101").Value = 4
$ws.Range("I102").Value = 7

# --- Selection moved by the author while reviewing the sheet ---
$ws.Range("I86").Select()
